$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows ("Fifth"/"sixth") that are no longer part of the data set
$ws.Rows("4:5").Delete()

# Relabel existing header (keep original styling/shared-string slots where possible)
$ws.Range("B1").Value = "DateTimeCol"

# New columns: parsed Date-only (C) and Time-only (D) values derived from the DateTimeCol text
$ws.Range("C1").Value = "DateCol"
$ws.Range("D1").Value = "TimeCol"

# Relabel remaining row label
$ws.Range("A2").Value = "Third"

# Row 2 (23Dec2020:12:20:34 -> date 12/23/2020, time 12:20:34 PM)
# Apply the time format first so the new styles are appended in the same
# order the workbook expects (time style before date style).
$ws.Range("D2").Value = 0.51428240740740738
$ws.Range("D2").NumberFormat = "[`$-F400]h:mm:ss\ AM/PM"
$ws.Range("C2").Value = 44188
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# Row 3 (23Jan2020:12:24:34 -> date 1/23/2020, time 00:24:35)
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 43853

$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 0.017071759259259259

# Drop the leftover explicit row heights so rows use the sheet default again
$ws.Rows("1:3").AutoFit()

# Column widths to roughly match the authored layout (ColumnWidth uses a
# slightly different unit than the stored XML width, so back it off by the
# engine's fixed offset to land as close as possible to the target width).
$ws.Columns("D:D").ColumnWidth = 8.998697916666666

$ws.Range("A4").Select() | Out-Null
